# Update the "想去人数" (attendance) figures for two events that each
# appear in both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    if ($name -eq "展览") {
        $ws.Range("F3").Value = 7560
        $ws.Range("F6").Value = 461
        $ws.Range("F7").Value = 4188
    }
    else {
        $ws.Range("F4").Value = 7560
        $ws.Range("F8").Value = 461
        $ws.Range("F9").Value = 4188
    }
}
